$wb = $excel.ActiveWorkbook

# Add the new worksheet and move it to the end (after NullableClass)
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "RowHeader"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

$ws = $wb.Worksheets.Item("RowHeader")

# Row headers (column A), written top to bottom first
$ws.Range("A1").Value = "Name"
$ws.Range("A2").Value = "Age"
$ws.Range("A3").Value = "Username"
$ws.Range("A4").Value = "Email"
$ws.Range("A5").Value = "Creation Date"

# Robert's column (B), written top to bottom
$ws.Range("B1").Value = "Robert"
$ws.Range("B2").Value = 42
$ws.Range("B3").Value = "Robert42"
$ws.Range("B4").Value = "robert@gmail.com"
$ws.Range("B5").Value = 27954
$ws.Range("B5").NumberFormat = "mm-dd-yy"

# Adam's column (C), written top to bottom
$ws.Range("C1").Value = "Adam"
$ws.Range("C2").Value = 28
$ws.Range("C3").Value = "Adam28"
$ws.Range("C4").Value = "adam@gmail.com"
$ws.Range("C5").Value = 32153
$ws.Range("C5").NumberFormat = "mm-dd-yy"

# Hyperlinks on the email cells
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:robert@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:adam@gmail.com")

# Column widths
$ws.Columns.Item(1).ColumnWidth = 13.26953125
$ws.Columns.Item(2).ColumnWidth = 19.08984375
$ws.Columns.Item(3).ColumnWidth = 26.453125

# Row height for row 4 (slightly taller, matches ht="16")
$ws.Rows.Item(4).RowHeight = 16

# Update the selection on sheet "Projects" (sheet1) to B5
$sheet1 = $wb.Worksheets.Item("Projects")
$sheet1.Range("B5").Select() | Out-Null

# Selection on the new sheet, then activate it so it is the visible tab
$ws.Range("D8").Select() | Out-Null
$ws.Activate() | Out-Null
